$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.155.48'
$ws.Range('E2').Value = '  +1.56%  '

$ws.Range('D3').Value = '3.118.77'
$ws.Range('E3').Value = '  -1.88%  '

$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.19'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '615.57'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.72%  '

$ws.Range('E7').Value = '  +0.49%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.399'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.37%  '

$ws.Range('E9').Value = '  -0.04%  '

$ws.Range('D10').Value = '3.117.38'
$ws.Range('E10').Value = '  +30.50%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.755'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('E12').Value = '  -0.83%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.78%  '

$ws.Range('D14').Value = '93.202.91'
$ws.Range('E14').Value = '  +2.04%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '34.54'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.70%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.47'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.46%  '

$ws.Range('D17').Value = '3.714.77'
$ws.Range('E17').Value = '  -0.89%  '

$ws.Range('D18').Value = '3.125.18'
$ws.Range('E18').Value = '  -0.98%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.79'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.75%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.73%  '

$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.80'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.12%  '

$ws.Range('B22').Value = 'PEPE'
$ws.Range('C22').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0000205'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '448.46'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.19%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.34'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.66%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.81'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.09%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '87.18'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.54%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.78'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.58%  '

$ws.Range('D28').Value = '3.294.05'
$ws.Range('E28').Value = '  -0.63%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.136'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +12.26%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.231'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.23%  '

$ws.Range('E32').Value = '  -0.44%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '9.23'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.67%  '

$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '8.12'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.93%  '

$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.160'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -6.02%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '26.17'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.73%  '

$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.879'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.29%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '493.22'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.03%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.90'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.49%  '

$ws.Range('E40').Value = '  +1.09%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.30'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.90%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.433'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.43%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.41'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.08%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.09'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +4.12%  '

$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '164.11'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.90%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.92'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.70%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.689'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.87%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.39'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.43%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0334'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +4.92%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.43'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.02%  '
